$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.4128796458244324
$ws.Range("B1").Value = 0.3068462014198303
$ws.Range("C1").Value = 0.3703226745128632
$ws.Range("D1").Value = 3.905551671981812
$ws.Range("E1").Value = 1.653122663497925
